# "business questions code is updated"
#
# Adds a new "TC004" business-question row to Sheet1 of the data-driven
# test sheet (a clone of the TC003 row, with the test-case id updated and
# the HomeOwnerDisc flag flipped to "Yes"), wires up the matching mailto
# hyperlink for the new row's e-mail cell, and leaves Sheet1 active with
# the newly-added row selected - mirroring what happens when a tester adds
# one more scripted test case to the data sheet and saves.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$srcRow = 4   # TC003
$dstRow = 5   # new TC004

# --- Clone the whole TC003 row (values + number formats/styles) into the
#     new row, then adjust the two cells that differ for TC004.
$ws1.Range("A" + $srcRow + ":AC" + $srcRow).Copy() | Out-Null
$ws1.Range("A" + $dstRow + ":AC" + $dstRow).PasteSpecial(-4163) | Out-Null # xlPasteValues
$ws1.Range("A" + $srcRow + ":AC" + $srcRow).Copy() | Out-Null
$ws1.Range("A" + $dstRow + ":AC" + $dstRow).PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

# New test case id.
$ws1.Cells.Item($dstRow, 1).Value2 = "TC004"

# HomeOwnerDisc (column AA) is "Yes" for TC004 (it was "No" for TC003), and
# picks up the text-number-format style already used elsewhere on the row.
$ws1.Cells.Item($dstRow, 27).Value2 = "Yes"
$ws1.Range("G" + $dstRow).Copy() | Out-Null
$ws1.Range("AA" + $dstRow).PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

# --- Hyperlink for the new row's e-mail cell (same mailto pattern as the
#     other rows). Adding it restyles the cell, so restore its original
#     (non-hyperlink-default) look straight after.
$email = $ws1.Cells.Item($dstRow, 2).Value2
$ws1.Hyperlinks.Add($ws1.Cells.Item($dstRow, 2), "mailto:" + $email) | Out-Null
$ws1.Range("B" + $srcRow).Copy() | Out-Null
$ws1.Range("B" + $dstRow).PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

# --- Make Sheet1 the active sheet with the new row's last cell selected.
$ws1.Activate()
$ws1.Range("AA" + $dstRow).Select() | Out-Null
